$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 171, shifting existing rows 171:197 down to 172:198
$ws.Rows("171:171").Insert()

# Populate the newly inserted row with the new data record
$ws.Cells.Item(171, 1).Value2 = 9
$ws.Cells.Item(171, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(171, 3).Value2 = "Metropolitana"
$ws.Cells.Item(171, 4).Value2 = 44637
$ws.Cells.Item(171, 5).Value2 = 13
$ws.Cells.Item(171, 6).Value2 = 100112026
$ws.Cells.Item(171, 7).Value2 = "Haba"
$ws.Cells.Item(171, 8).Value2 = "Sin especificar"
$ws.Cells.Item(171, 9).Value2 = "Primera"
$ws.Cells.Item(171, 10).Value2 = 20
$ws.Cells.Item(171, 11).Value2 = 20000
$ws.Cells.Item(171, 12).Value2 = 20000
$ws.Cells.Item(171, 13).Value2 = 20000
$ws.Cells.Item(171, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(171, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(171, 16).Value2 = 800
$ws.Cells.Item(171, 17).Value2 = 25
$ws.Cells.Item(171, 18).Value2 = "Hortaliza"
